$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03958166666666667
$ws.Range("H2").Value = 0.118745
$ws.Range("I2").Value = 0.2870281964201545
$ws.Range("J2").Value = 0.2870281964201545
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.51469433333333
$ws.Range("N2").Value = 67.544083
$ws.Range("O2").Value = 0.1309757462958079
$ws.Range("P2").Value = 0.1309757462958079
$ws.Range("Q2").Value = 0.891169126203889
$ws.Range("R2").Value = 8.020522135835
$ws.Range("S2").Value = 0.03759373223406948
$ws.Range("T2").Value = 0.03759373223406948

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03958166666666667
$ws.Range("H3").Value = 0.118745
$ws.Range("I3").Value = 0.2870281964201545
$ws.Range("J3").Value = 0.2870281964201545
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.26507466666667
$ws.Range("N3").Value = 57.795224
$ws.Range("O3").Value = 0.1120715873177727
$ws.Range("P3").Value = 0.1120715873177728
$ws.Range("Q3").Value = 0.7625437637644444
$ws.Range("R3").Value = 6.86289387388
$ws.Range("S3").Value = 0.03216770557776417
$ws.Range("T3").Value = 0.03216770557776417

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03958166666666667
$ws.Range("H4").Value = 0.118745
$ws.Range("I4").Value = 0.2870281964201545
$ws.Range("J4").Value = 0.2870281964201545
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 46.79779933333333
$ws.Range("N4").Value = 140.393398
$ws.Range("O4").Value = 0.272238947681833
$ws.Range("P4").Value = 0.272238947681833
$ws.Range("Q4").Value = 1.852334893945556
$ws.Range("R4").Value = 16.67101404551
$ws.Range("S4").Value = 0.07814025414843732
$ws.Range("T4").Value = 0.07814025414843734

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03958166666666667
$ws.Range("H5").Value = 0.118745
$ws.Range("I5").Value = 0.2870281964201545
$ws.Range("J5").Value = 0.2870281964201545
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.237739333333333
$ws.Range("N5").Value = 18.713218
$ws.Range("O5").Value = 0.03628708221778873
$ws.Range("P5").Value = 0.03628708221778874
$ws.Range("Q5").Value = 0.2469001190455555
$ws.Range("R5").Value = 2.22210107141
$ws.Range("S5").Value = 0.01041541576232176
$ws.Range("T5").Value = 0.01041541576232176

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03958166666666667
$ws.Range("H6").Value = 0.118745
$ws.Range("I6").Value = 0.2870281964201545
$ws.Range("J6").Value = 0.2870281964201545
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 35.481449
$ws.Range("N6").Value = 106.444347
$ws.Range("O6").Value = 0.2064078327526475
$ws.Range("P6").Value = 0.2064078327526475
$ws.Range("Q6").Value = 1.404414887168333
$ws.Range("R6").Value = 12.639733984515
$ws.Range("S6").Value = 0.0592448679619853
$ws.Range("T6").Value = 0.05924486796198531

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03958166666666667
$ws.Range("H7").Value = 0.118745
$ws.Range("I7").Value = 0.2870281964201545
$ws.Range("J7").Value = 0.2870281964201545
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 41.602965
$ws.Range("N7").Value = 124.808895
$ws.Range("O7").Value = 0.2420188037341499
$ws.Range("P7").Value = 0.2420188037341499
$ws.Range("Q7").Value = 1.646714692975
$ws.Range("R7").Value = 14.820432236775
$ws.Range("S7").Value = 0.0694662207355764
$ws.Range("T7").Value = 0.0694662207355764

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.09832
$ws.Range("H8").Value = 0.29496
$ws.Range("I8").Value = 0.7129718035798456
$ws.Range("J8").Value = 0.7129718035798456
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.51469433333333
$ws.Range("N8").Value = 67.544083
$ws.Range("O8").Value = 0.1309757462958079
$ws.Range("P8").Value = 0.1309757462958079
$ws.Range("Q8").Value = 2.213644746853334
$ws.Range("R8").Value = 19.92280272168
$ws.Range("S8").Value = 0.09338201406173846
$ws.Range("T8").Value = 0.09338201406173846

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.09832
$ws.Range("H9").Value = 0.29496
$ws.Range("I9").Value = 0.7129718035798456
$ws.Range("J9").Value = 0.7129718035798456
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.26507466666667
$ws.Range("N9").Value = 57.795224
$ws.Range("O9").Value = 0.1120715873177727
$ws.Range("P9").Value = 0.1120715873177728
$ws.Range("Q9").Value = 1.894142141226667
$ws.Range("R9").Value = 17.04727927104
$ws.Range("S9").Value = 0.07990388174000858
$ws.Range("T9").Value = 0.07990388174000859

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.09832
$ws.Range("H10").Value = 0.29496
$ws.Range("I10").Value = 0.7129718035798456
$ws.Range("J10").Value = 0.7129718035798456
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 46.79779933333333
$ws.Range("N10").Value = 140.393398
$ws.Range("O10").Value = 0.272238947681833
$ws.Range("P10").Value = 0.272238947681833
$ws.Range("Q10").Value = 4.601159630453333
$ws.Range("R10").Value = 41.41043667408
$ws.Range("S10").Value = 0.1940986935333957
$ws.Range("T10").Value = 0.1940986935333957

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.09832
$ws.Range("H11").Value = 0.29496
$ws.Range("I11").Value = 0.7129718035798456
$ws.Range("J11").Value = 0.7129718035798456
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.237739333333333
$ws.Range("N11").Value = 18.713218
$ws.Range("O11").Value = 0.03628708221778873
$ws.Range("P11").Value = 0.03628708221778874
$ws.Range("Q11").Value = 0.6132945312533333
$ws.Range("R11").Value = 5.519650781279999
$ws.Range("S11").Value = 0.02587166645546697
$ws.Range("T11").Value = 0.02587166645546698

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.09832
$ws.Range("H12").Value = 0.29496
$ws.Range("I12").Value = 0.7129718035798456
$ws.Range("J12").Value = 0.7129718035798456
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 35.481449
$ws.Range("N12").Value = 106.444347
$ws.Range("O12").Value = 0.2064078327526475
$ws.Range("P12").Value = 0.2064078327526475
$ws.Range("Q12").Value = 3.48853606568
$ws.Range("R12").Value = 31.39682459112
$ws.Range("S12").Value = 0.1471629647906622
$ws.Range("T12").Value = 0.1471629647906622

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.09832
$ws.Range("H13").Value = 0.29496
$ws.Range("I13").Value = 0.7129718035798456
$ws.Range("J13").Value = 0.7129718035798456
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 41.602965
$ws.Range("N13").Value = 124.808895
$ws.Range("O13").Value = 0.2420188037341499
$ws.Range("P13").Value = 0.2420188037341499
$ws.Range("Q13").Value = 4.090403518800001
$ws.Range("R13").Value = 36.8136316692
$ws.Range("S13").Value = 0.1725525829985735
$ws.Range("T13").Value = 0.1725525829985735
